$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 1).Value = 0
$ws.Cells.Item(2, 2).Value = "walkingToRunning"
$ws.Cells.Item(2, 3).Value = 18.70621109008789
$ws.Cells.Item(2, 4).Value = -5.636280059814453
$ws.Cells.Item(2, 5).Value = 6.914060115814209
$ws.Cells.Item(2, 6).Value = -1.059356399572612
$ws.Cells.Item(2, 7).Value = -0.06877905879515356
$ws.Cells.Item(2, 8).Value = 2.407926618075761

$ws.Cells.Item(3, 1).Value = 100
$ws.Cells.Item(3, 2).Value = "walkingToRunning"
$ws.Cells.Item(3, 3).Value = 19.93985939025879
$ws.Cells.Item(3, 4).Value = -3.599599838256836
$ws.Cells.Item(3, 5).Value = 8.073759078979492
$ws.Cells.Item(3, 6).Value = -1.104726179701384
$ws.Cells.Item(3, 7).Value = 1.440624993355549
$ws.Cells.Item(3, 8).Value = 1.773343285576245

$ws.Cells.Item(4, 1).Value = 200
$ws.Cells.Item(4, 2).Value = "walkingToRunning"
$ws.Cells.Item(4, 3).Value = 27.86785316467285
$ws.Cells.Item(4, 4).Value = -18.30047225952148
$ws.Cells.Item(4, 5).Value = 6.377731800079346
$ws.Cells.Item(4, 6).Value = -0.5441059966882074
$ws.Cells.Item(4, 7).Value = 0.8777113159497576
$ws.Cells.Item(4, 8).Value = 1.207277496655782

$ws.Cells.Item(5, 1).Value = 300
$ws.Cells.Item(5, 2).Value = "walkingToRunning"
$ws.Cells.Item(5, 3).Value = -0.164954662322998
$ws.Cells.Item(5, 4).Value = -10.53660011291504
$ws.Cells.Item(5, 5).Value = 4.150550365447998
$ws.Cells.Item(5, 6).Value = 0.1515757088127038
$ws.Cells.Item(5, 7).Value = 0.07700501178783151
$ws.Cells.Item(5, 8).Value = 0.8280433637848317

$ws.Cells.Item(6, 1).Value = 400
$ws.Cells.Item(6, 2).Value = "walkingToRunning"
$ws.Cells.Item(6, 3).Value = 5.472853660583496
$ws.Cells.Item(6, 4).Value = -46.38559722900391
$ws.Cells.Item(6, 5).Value = -0.9001345634460449
$ws.Cells.Item(6, 6).Value = -0.8091414912802259
$ws.Cells.Item(6, 7).Value = 0.04339970233010434
$ws.Cells.Item(6, 8).Value = 1.094123736756746

$ws.Cells.Item(7, 1).Value = 500
$ws.Cells.Item(7, 2).Value = "walkingToRunning"
$ws.Cells.Item(7, 3).Value = 3.861607551574707
$ws.Cells.Item(7, 4).Value = -19.79164886474609
$ws.Cells.Item(7, 5).Value = 0.622889518737793
$ws.Cells.Item(7, 6).Value = -1.296450522428008
$ws.Cells.Item(7, 7).Value = 0.4430659739046126
$ws.Cells.Item(7, 8).Value = 1.507120371190581

$ws.Cells.Item(8, 1).Value = 600
$ws.Cells.Item(8, 2).Value = "walkingToRunning"
$ws.Cells.Item(8, 3).Value = 14.61344242095947
$ws.Cells.Item(8, 4).Value = -10.30790042877197
$ws.Cells.Item(8, 5).Value = 17.13724899291992
$ws.Cells.Item(8, 6).Value = 1.030188943722537
$ws.Cells.Item(8, 7).Value = 1.112867043969416
$ws.Cells.Item(8, 8).Value = -0.4273038459280109

$ws.Cells.Item(9, 1).Value = 700
$ws.Cells.Item(9, 2).Value = "walkingToRunning"
$ws.Cells.Item(9, 3).Value = 2.93750786781311
$ws.Cells.Item(9, 4).Value = -46.01893997192383
$ws.Cells.Item(9, 5).Value = 7.825448036193848
$ws.Cells.Item(9, 6).Value = -0.01950983642064261
$ws.Cells.Item(9, 7).Value = -2.121916278463881
$ws.Cells.Item(9, 8).Value = -2.743124836780999

$ws.Cells.Item(10, 1).Value = 800
$ws.Cells.Item(10, 2).Value = "walkingToRunning"
$ws.Cells.Item(10, 3).Value = -20.74547958374023
$ws.Cells.Item(10, 4).Value = -0.2439025640487671
$ws.Cells.Item(10, 5).Value = 8.103152275085449
$ws.Cells.Item(10, 6).Value = 2.686820753651581
$ws.Cells.Item(10, 7).Value = -0.1019303134230194
$ws.Cells.Item(10, 8).Value = -4.739581848102853

$ws.Cells.Item(11, 1).Value = 900
$ws.Cells.Item(11, 2).Value = "walkingToRunning"
$ws.Cells.Item(11, 3).Value = -20.94050025939941
$ws.Cells.Item(11, 4).Value = -7.066665649414063
$ws.Cells.Item(11, 5).Value = 6.037558555603027
$ws.Cells.Item(11, 6).Value = 2.067020179144016
$ws.Cells.Item(11, 7).Value = 4.138878788453006
$ws.Cells.Item(11, 8).Value = -7.458814548013004

$ws.Cells.Item(12, 1).Value = 1000
$ws.Cells.Item(12, 2).Value = "walkingToRunning"
$ws.Cells.Item(12, 3).Value = 29.87363624572754
$ws.Cells.Item(12, 4).Value = -61.11213684082031
$ws.Cells.Item(12, 5).Value = 14.88786697387695
$ws.Cells.Item(12, 6).Value = -1.318529965447579
$ws.Cells.Item(12, 7).Value = -0.9071916752173932
$ws.Cells.Item(12, 8).Value = 2.900026571555154

$ws.Cells.Item(13, 1).Value = 1100
$ws.Cells.Item(13, 2).Value = "walkingToRunning"
$ws.Cells.Item(13, 3).Value = 4.961847305297852
$ws.Cells.Item(13, 4).Value = -15.47994041442871
$ws.Cells.Item(13, 5).Value = 18.70905494689941
$ws.Cells.Item(13, 6).Value = -5.118491420328921
$ws.Cells.Item(13, 7).Value = -0.9277605919238714
$ws.Cells.Item(13, 8).Value = 7.022751289638649

$ws.Cells.Item(14, 1).Value = 1200
$ws.Cells.Item(14, 2).Value = "walkingToRunning"
$ws.Cells.Item(14, 3).Value = -0.8285019397735596
$ws.Cells.Item(14, 4).Value = -6.118541240692139
$ws.Cells.Item(14, 5).Value = 8.952471733093262
$ws.Cells.Item(14, 6).Value = -4.701338476170607
$ws.Cells.Item(14, 7).Value = 4.033473036979721
$ws.Cells.Item(14, 8).Value = 6.285576557201113

$ws.Cells.Item(15, 1).Value = 1300
$ws.Cells.Item(15, 2).Value = "walkingToRunning"
$ws.Cells.Item(15, 3).Value = 0.8424484729766846
$ws.Cells.Item(15, 4).Value = -38.40699005126953
$ws.Cells.Item(15, 5).Value = 14.60053634643555
$ws.Cells.Item(15, 6).Value = 0.6423424071953132
$ws.Cells.Item(15, 7).Value = -1.238872872024314
$ws.Cells.Item(15, 8).Value = 0.4076568924012518

$ws.Cells.Item(16, 1).Value = 1400
$ws.Cells.Item(16, 2).Value = "walkingToRunning"
$ws.Cells.Item(16, 3).Value = 1.708237409591675
$ws.Cells.Item(16, 4).Value = 32.82785034179688
$ws.Cells.Item(16, 5).Value = 13.28276348114014
$ws.Cells.Item(16, 6).Value = 2.748690418858344
$ws.Cells.Item(16, 7).Value = -5.918399258389504
$ws.Cells.Item(16, 8).Value = -5.032549650942702

$ws.Cells.Item(17, 1).Value = 1500
$ws.Cells.Item(17, 2).Value = "walkingToRunning"
$ws.Cells.Item(17, 3).Value = -19.85161781311035
$ws.Cells.Item(17, 4).Value = -6.960978984832764
$ws.Cells.Item(17, 5).Value = 4.844282150268555
$ws.Cells.Item(17, 6).Value = 3.90800370544684
$ws.Cells.Item(17, 7).Value = 1.824196169285198
$ws.Cells.Item(17, 8).Value = -6.420193484572113

$ws.Cells.Item(18, 1).Value = 1600
$ws.Cells.Item(18, 2).Value = "walkingToRunning"
$ws.Cells.Item(18, 3).Value = -47.64518356323242
$ws.Cells.Item(18, 4).Value = -56.76200866699219
$ws.Cells.Item(18, 5).Value = 40.71841812133789
$ws.Cells.Item(18, 6).Value = -0.09310042271849372
$ws.Cells.Item(18, 7).Value = 3.190837242564208
$ws.Cells.Item(18, 8).Value = -1.485911830526877

$ws.Cells.Item(19, 1).Value = 1700
$ws.Cells.Item(19, 2).Value = "walkingToRunning"
$ws.Cells.Item(19, 3).Value = -2.008986234664917
$ws.Cells.Item(19, 4).Value = -3.396074771881104
$ws.Cells.Item(19, 5).Value = 5.574520111083984
$ws.Cells.Item(19, 6).Value = -5.086900020557672
$ws.Cells.Item(19, 7).Value = -0.347581442587992
$ws.Cells.Item(19, 8).Value = 4.233577517212423

$ws.Cells.Item(20, 1).Value = 1800
$ws.Cells.Item(20, 2).Value = "walkingToRunning"
$ws.Cells.Item(20, 3).Value = 1.793292045593261
$ws.Cells.Item(20, 4).Value = 1.589181900024414
$ws.Cells.Item(20, 5).Value = 13.23852920532227
$ws.Cells.Item(20, 6).Value = -4.690330617414796
$ws.Cells.Item(20, 7).Value = 3.666637846680965
$ws.Cells.Item(20, 8).Value = 4.422437811158368

$ws.Cells.Item(21, 1).Value = 1900
$ws.Cells.Item(21, 2).Value = "walkingToRunning"
$ws.Cells.Item(21, 3).Value = 0.5213950872421265
$ws.Cells.Item(21, 4).Value = -13.69121932983398
$ws.Cells.Item(21, 5).Value = 13.58244514465332
$ws.Cells.Item(21, 6).Value = -3.853545032563759
$ws.Cells.Item(21, 7).Value = 7.179606515853127
$ws.Cells.Item(21, 8).Value = 0.02150726318360263

$ws.Cells.Item(22, 1).Value = 2000
$ws.Cells.Item(22, 2).Value = "walkingToRunning"
$ws.Cells.Item(22, 3).Value = -25.09261322021484
$ws.Cells.Item(22, 4).Value = 15.60748481750488
$ws.Cells.Item(22, 5).Value = 0.5675735473632812
$ws.Cells.Item(22, 6).Value = -0.654346335781087
$ws.Cells.Item(22, 7).Value = 7.300580988816252
$ws.Cells.Item(22, 8).Value = -5.563441323452309

$ws.Cells.Item(23, 1).Value = 2100
$ws.Cells.Item(23, 2).Value = "walkingToRunning"
$ws.Cells.Item(23, 3).Value = -30.59898567199707
$ws.Cells.Item(23, 4).Value = -12.55906105041504
$ws.Cells.Item(23, 5).Value = 3.974555969238281
$ws.Cells.Item(23, 6).Value = 3.714518678644319
$ws.Cells.Item(23, 7).Value = -2.467024099305689
$ws.Cells.Item(23, 8).Value = -5.610007460651495

$ws.Cells.Item(24, 1).Value = 2200
$ws.Cells.Item(24, 2).Value = "walkingToRunning"
$ws.Cells.Item(24, 3).Value = -38.86380767822266
$ws.Cells.Item(24, 4).Value = -84.71040344238281
$ws.Cells.Item(24, 5).Value = 66.18233489990234
$ws.Cells.Item(24, 6).Value = 2.151245980966297
$ws.Cells.Item(24, 7).Value = 0.7463607152954532
$ws.Cells.Item(24, 8).Value = -6.372515623686779

$ws.Cells.Item(25, 1).Value = 2300
$ws.Cells.Item(25, 2).Value = "walkingToRunning"
$ws.Cells.Item(25, 3).Value = -8.86505126953125
$ws.Cells.Item(25, 4).Value = -2.532943964004517
$ws.Cells.Item(25, 5).Value = 5.462150573730469
$ws.Cells.Item(25, 6).Value = -1.10870781817722
$ws.Cells.Item(25, 7).Value = -0.5732807598478937
$ws.Cells.Item(25, 8).Value = 3.703908173764338

$ws.Cells.Item(26, 1).Value = 2400
$ws.Cells.Item(26, 2).Value = "walkingToRunning"
$ws.Cells.Item(26, 3).Value = 7.391067981719971
$ws.Cells.Item(26, 4).Value = -1.471791982650757
$ws.Cells.Item(26, 5).Value = 18.45427322387696
$ws.Cells.Item(26, 6).Value = -3.764231722211575
$ws.Cells.Item(26, 7).Value = 3.633465962331838
$ws.Cells.Item(26, 8).Value = 3.590699133977193

$ws.Cells.Item(27, 1).Value = 2500
$ws.Cells.Item(27, 2).Value = "walkingToRunning"
$ws.Cells.Item(27, 3).Value = -4.55918025970459
$ws.Cells.Item(27, 4).Value = -21.72416114807129
$ws.Cells.Item(27, 5).Value = -0.08992767333984369
$ws.Cells.Item(27, 6).Value = -3.505413231302479
$ws.Cells.Item(27, 7).Value = 9.71845419680489
$ws.Cells.Item(27, 8).Value = 0.237572531231097

$ws.Cells.Item(28, 1).Value = 2600
$ws.Cells.Item(28, 2).Value = "walkingToRunning"
$ws.Cells.Item(28, 3).Value = -17.91468048095703
$ws.Cells.Item(28, 4).Value = 18.05105400085449
$ws.Cells.Item(28, 5).Value = -10.55736446380615
$ws.Cells.Item(28, 6).Value = -2.849029684327336
$ws.Cells.Item(28, 7).Value = 8.877369792083568
$ws.Cells.Item(28, 8).Value = -7.541428113895723

$ws.Cells.Item(29, 1).Value = 2700
$ws.Cells.Item(29, 2).Value = "walkingToRunning"
$ws.Cells.Item(29, 3).Value = -62.19514083862305
$ws.Cells.Item(29, 4).Value = -14.345703125
$ws.Cells.Item(29, 5).Value = 3.762966632843018
$ws.Cells.Item(29, 6).Value = 3.086845004493447
$ws.Cells.Item(29, 7).Value = -5.60899220659439
$ws.Cells.Item(29, 8).Value = -5.2191632692931

$ws.Cells.Item(30, 1).Value = 2800
$ws.Cells.Item(30, 2).Value = "walkingToRunning"
$ws.Cells.Item(30, 3).Value = -47.85998153686523
$ws.Cells.Item(30, 4).Value = -73.05361175537109
$ws.Cells.Item(30, 5).Value = 52.46365737915039
$ws.Cells.Item(30, 6).Value = 3.54143344769713
$ws.Cells.Item(30, 7).Value = -2.492278899325752
$ws.Cells.Item(30, 8).Value = -5.184806862815476

$ws.Cells.Item(31, 1).Value = 2900
$ws.Cells.Item(31, 2).Value = "walkingToRunning"
$ws.Cells.Item(31, 3).Value = -2.905624389648437
$ws.Cells.Item(31, 4).Value = 1.127065658569336
$ws.Cells.Item(31, 5).Value = 7.177680492401123
$ws.Cells.Item(31, 6).Value = -0.4025364195714203
$ws.Cells.Item(31, 7).Value = -2.233978587095854
$ws.Cells.Item(31, 8).Value = 0.2824505248356344
